# Adjusting figure elements in experimental panel:
# add two small green "$$" / "$" callout textboxes near the
# exploration/exploitation picture group on slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# PowerPoint COM Shapes.AddTextbox / shape geometry is expressed in points;
# the target OOXML offsets/extents are in EMU (1 pt = 12700 EMU), so convert.
$emuPerPt = 12700

# --- consume the two shape-id slots ("19" and "21") that the original
#     authoring session burned through (e.g. via an add+undo) before the
#     two kept textboxes landed on ids 20 and 22 ---
$burn1 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$burn1.Delete()

# --- TextBox "TextBox 19" (ends up as shape id 20) ---
$tb1Left   = 1092707 / $emuPerPt
$tb1Top    = 1412227 / $emuPerPt
$tb1Width  = 439134 / $emuPerPt
$tb1Height = 369332 / $emuPerPt

$tb1 = $s.Shapes.AddTextbox(1, $tb1Left, $tb1Top, $tb1Width, $tb1Height)
$tb1.Name = "TextBox 19"
$tb1.Fill.Visible = 0
$tb1.TextFrame.WordWrap = 1
$tb1.TextFrame.AutoSize = 1

$tr1 = $tb1.TextFrame.TextRange
$tr1.Text = "`$`$"
$tr1.ParagraphFormat.Alignment = 2
$tr1.Font.Color.RGB = 5287936
$tr1.Font.NameFarEast = "Calibri"
$tr1.Font.NameComplexScript = "Calibri"

# --- consume the next shape-id slot ("21") ---
$burn2 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$burn2.Delete()

# --- TextBox "TextBox 21" (ends up as shape id 22) ---
$tb2Left   = 2330391 / $emuPerPt
$tb2Top    = 1985407 / $emuPerPt
$tb2Width  = 439134 / $emuPerPt
$tb2Height = 369332 / $emuPerPt

$tb2 = $s.Shapes.AddTextbox(1, $tb2Left, $tb2Top, $tb2Width, $tb2Height)
$tb2.Name = "TextBox 21"
$tb2.Fill.Visible = 0
$tb2.TextFrame.WordWrap = 1
$tb2.TextFrame.AutoSize = 1

$tr2 = $tb2.TextFrame.TextRange
$tr2.Text = "`$"
$tr2.ParagraphFormat.Alignment = 2
$tr2.Font.Color.RGB = 5287936
$tr2.Font.NameFarEast = "Calibri"
$tr2.Font.NameComplexScript = "Calibri"
